$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect so we can update the cells below.
$ws.Unprotect()

# --- Update the "as of" date in the confidential disclaimer text (A9) ---
$nl = [char]10
$disclaimerText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + $nl + "Model holdings provided as of 2021-04-08 for illustrative purposes only and are subject to change."
$ws.Range("A9").Value = $disclaimerText

# --- Update the Weight (D) and Percent Change (E) figures for rows 2-6 ---
$ws.Range("D2").Value = 0.2466256200330046
$ws.Range("E2").Value = 0.001261670451678132

$ws.Range("D3").Value = 0.2471783317474865
$ws.Range("E3").Value = 0.0002870264064291916

$ws.Range("D4").Value = 0.2553790122870609
$ws.Range("E4").Value = 0.01444877907816777

$ws.Range("D5").Value = 0.2508170359324481
$ws.Range("E5").Value = 0.0001299883010528191

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.004104625175843424

# Restore sheet protection to match the original protected state.
$ws.Protect("D382", $false, $true, $false, $false)
